# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) for the affected leve rows across the profession
# sheets. Values below reflect the latest Universalis pull.

$wb = $excel.ActiveWorkbook

function Set-Cells($SheetName, $Row, $Values) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $addr = "$col$Row"
        $val = $Values[$col]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}

# ALC
Set-Cells "ALC" 28  @{ H=1199.5358; J=1098.1111; L=1098.1111; N=-2068.1111 }
Set-Cells "ALC" 32  @{ H=12479.9; I=11666.667; J=12828.429; K=11666.667; L=12828.429; M=-11340.667; N=-13480.429 }
Set-Cells "ALC" 74  @{ H=55927.15; I=71526.266; K=71526.266; M=-70590.266 }
Set-Cells "ALC" 77  @{ H=55927.15; I=71526.266; K=357631.33; M=-352951.33 }
Set-Cells "ALC" 98  @{ H=491.3; I=491.3; J=0; K=491.3; L=0; M=1006.7; N=$null }
Set-Cells "ALC" 122 @{ H=491.3; I=491.3; J=0; K=1473.9; L=0; M=976.0999999999999; N=$null }
Set-Cells "ALC" 137 @{ H=9011.143; I=2669.25; J=11547.9; K=8007.75; L=34643.7; M=-5457.75; N=-39743.7 }
Set-Cells "ALC" 138 @{ H=3158.078; J=3798.3333; L=11394.9999; N=-21674.9999 }
Set-Cells "ALC" 141 @{ H=151476.25; I=1950; J=301002.5; K=5850; L=903007.5; M=-670; N=-913367.5 }

# ARM
Set-Cells "ARM" 2   @{ H=1617.6072; I=1719.6; K=1719.6; M=-1606.6 }
Set-Cells "ARM" 74  @{ H=16071.429; J=34396.668; L=34396.668; N=-36144.668 }
Set-Cells "ARM" 77  @{ H=16071.429; J=34396.668; L=171983.34; N=-180719.34 }
Set-Cells "ARM" 110 @{ H=7715.778; I=9426; K=9426; M=-7381 }
Set-Cells "ARM" 116 @{ H=1617.6072; I=1719.6; K=1719.6; M=574.4000000000001 }

# BSM
Set-Cells "BSM" 3  @{ H=1617.6072; I=1719.6; K=1719.6; M=-1605.6 }
Set-Cells "BSM" 94 @{ H=2192.3 }

# CRP
Set-Cells "CRP" 31  @{ H=6704064.5; I=11152907; J=30799.9; K=11152907; L=30799.9; M=-11152612; N=-31389.9 }
Set-Cells "CRP" 34  @{ H=6704064.5; I=11152907; J=30799.9; K=11152907; L=30799.9; M=-11152705; N=-31203.9 }
Set-Cells "CRP" 99  @{ H=3806.8572; I=3806.8572; K=3806.8572; M=-2308.8572 }
Set-Cells "CRP" 126 @{ H=3806.8572; I=3806.8572; K=11420.5716; M=-8950.571599999999 }
Set-Cells "CRP" 132 @{ H=5666.6665; I=5000; J=7000; K=15000; L=21000; M=-12470; N=-26060 }
Set-Cells "CRP" 134 @{ H=11751.296; I=2614.25; K=7842.75; M=-5307.75 }

# CUL
Set-Cells "CUL" 12  @{ H=108.875; I=114.75; J=103; K=344.25; L=309; M=-171.25; N=-655 }
Set-Cells "CUL" 80  @{ H=3586; J=3673; L=11019; N=-12891 }
Set-Cells "CUL" 83  @{ H=3586; J=3673; L=33057; N=-42417 }
Set-Cells "CUL" 113 @{ H=554.44446; I=548.75; K=1646.25; M=523.75 }
Set-Cells "CUL" 124 @{ H=18576.154; I=5676; K=17028; M=-12118 }
Set-Cells "CUL" 131 @{ H=3988.182; I=345.77777; K=1037.33331; M=4002.66669 }
Set-Cells "CUL" 140 @{ H=8157.55; I=8803; K=26409; M=-21229 }

# GSM
Set-Cells "GSM" 113 @{ H=3899.2593; J=4419.5; L=4419.5; N=-8759.5 }

# LTW
Set-Cells "LTW" 46  @{ H=4199.8887; I=1800; K=1800; M=-1612 }
Set-Cells "LTW" 132 @{ H=6642.143; I=6000; J=7123.75; K=18000; L=21371.25; M=-15470; N=-26431.25 }
Set-Cells "LTW" 136 @{ H=7043.923; J=10000; L=30000; N=-35100 }

# WVR
Set-Cells "WVR" 42  @{ H=1000000000; J=1000000000; L=1000000000; N=-1000000756 }
Set-Cells "WVR" 132 @{ H=2829.2; I=2832.5557; K=8497.667099999999; M=-5967.667099999999 }
Set-Cells "WVR" 140 @{ H=51967.375; J=51967.375; L=51967.375; N=-62327.375 }
